$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- Change 1 -------------------------------------------------------------
# "Objetivo" paragraph: the colon-space run (": "), the lone "A" run and the
# "umentar..." run are re-unified into a single run carrying the full text
# ": Aumentar el alcance del gimnasio, generando ingresos a través de
# membresías y servicios en línea."
$pObjetivo = Get-ParagraphByText("Objetivo: A")
if ($pObjetivo -ne $null) {
    $rng1 = $pObjetivo.Range
    $rng1.Find.Execute(
        ": A", $false, $false, $false, $false, $false, $true, 1, $false,
        ": A", 2) | Out-Null
}

# --- Change 2 -------------------------------------------------------------
# "Canales de Marketing" paragraph: fix the stray double space in
# "...emails con  marketing." -> "...emails con marketing." and re-split the
# tail into four runs: "s ", "con ", "marketing", "."
$pMarketing = Get-ParagraphByText("Canales de Marketing")
if ($pMarketing -ne $null) {
    $rng2 = $pMarketing.Range
    $found2 = $rng2.Find.Execute(
        "s con  marketing.", $false, $false, $false, $false, $false, $true, 1,
        $false, $null, 0)

    if ($found2) {
        $base = $rng2.Start

        $rA = $d.Range($base, $base + 2)
        $rA.Text = "s "
        $rA.Bold = 1
        $rA.Bold = 0

        $rB = $d.Range($base + 2, $base + 7)
        $rB.Text = "con "
        $rB.Bold = 1
        $rB.Bold = 0

        $rC = $d.Range($base + 6, $base + 15)
        $rC.Text = "marketing"
        $rC.Bold = 1
        $rC.Bold = 0

        $rD = $d.Range($base + 15, $base + 16)
        $rD.Text = "."
        $rD.Bold = 1
        $rD.Bold = 0
    }
}
